# Policy.xlsx rule table simplification
#
# The "既有保戶(保單)" rule table used to fire two actions against two
# imported POJOs (Policy + Insured): set the Policy's name/amount and an
# "Insured" id condition. This edit drops the Policy-related columns
# entirely and leaves a single ACTION column that sets the Insured's
# total amount instead.
#
# Concretely (on sheet "policy"):
#   - Import (B2) no longer references com...pojo.Policy, only Insured.
#   - Columns C ("ACTION" / A223456789 literal / 姓名 / 陳小亮 / 陳老爺) and
#     D ("ACTION" / 25-wide 姓名 column) are removed outright; the old
#     "保單總額" (E) and "policy" (F) columns shift left to become the new
#     C and D columns (Excel keeps their merges/widths/row heights
#     attached as it shifts).
#   - The remaining ACTION cell (now C8) fires a different rule action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("policy")

# Drop the two Policy-specific columns (old C:D). Everything to the right
# (保單總額 / policy columns, their merged cells D10:D11, and column widths)
# shifts left automatically, landing exactly on the target layout
# (A1:D11, merges A5:C5 and D10:D11).
$ws.Range("C:D").Delete()

# Import list: drop the now-unused Policy POJO, keep only Insured.
$ws.Range("B2").Value = 'com.redhat.prudential_poc.pojo.Insured'

# ACTION cell for the existing-policyholder rule table: instead of firing
# a brand-new Policy and setting its id/name/amount, just set the total
# amount on the already-bound $insured.
$ws.Range("C8").Value = '$insured.setTotoalAmt($param);'

# Keep the active selection in sync with the new layout (old D8 selection
# now lands on C8 after the column shift).
$ws.Range("C8").Select()
